$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current column B (Item1Name),
# shifting Item1.. / Item2.. / Item3.. columns two places to the right.
$ws.Range("B1:C1").EntireColumn.Insert()

# New header cells for the inserted columns.
$ws.Range("B1").Value = "ItemCountMin"
$ws.Range("C1").Value = "ItemCountMax"

# Populate the new ItemCountMin / ItemCountMax columns for each data row.
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 3

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 3

# Match the authored column widths for the two new columns
# (ColumnWidth values are pre-compensated for the engine's
# characters -> stored-width pixel rounding so the saved XML
# 'width' attribute lands on 12 / 12.5 exactly).
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(3).ColumnWidth = 11.666666666666666

# Move the selection like in the authored workbook.
$ws.Range("C6").Select()
